$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cameras")

# Row 7 - new Sony a7ii camera entry
$ws.Range("A7").Value = "Sony a7ii"
$ws.Range("B7").Value = 35.799999999999997
$ws.Range("C7").Value = 23.9
$ws.Range("D7").Value = 6000
$ws.Range("E7").Value = 4000

$ws.Range("F7").Formula = "=B7/D7*1000"
$ws.Range("G7").Formula = "=C7/E7*1000"
$ws.Range("H7").Formula = "=AVERAGE(F7:G7)"
$ws.Range("I7").Formula = "=I6"
$ws.Range("J7").Formula = "=I7/H7"
$ws.Range("K7").Formula = '=_xlfn.CONCAT("- ",ROUND(B7*J7,0),"x",ROUND(C7*J7,0)," mm")'
$ws.Range("L7").Formula = '=_xlfn.CONCAT("- ",ROUND(B7*J7/25.4,0),"x",ROUND(C7*J7/25.4,0)," inch")'

# Copy styles from row 6 to row 7
$ws.Range("F6:H6").Copy()
$ws.Range("F7:H7").PasteSpecial(-4122)
$ws.Range("I6").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("K6:L6").Copy()
$ws.Range("K7:L7").PasteSpecial(-4122)

# Update selection and view
$ws.Range("G7").Select()
$excel.ActiveWindow.ScrollRow = 14
